# design_subject_list.xlsx — apply "fixed multiple reference" commit
#
# Every data row (2-9) had its 개설학과(B)/교과구분(D)/학기(F) columns padded
# with a single-space placeholder string, and course row 2's name cell
# (C2, "기초창의공학설계") was missing the Korean-charset font style that every
# other course-name cell already carries. This pass clears the stray
# placeholders (keeping the cell formatting) and normalizes C2's font to
# match its siblings, then nudges the active selection to where editing
# left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the leftover " " placeholder values in columns B, D, F for rows 2-9.
# ClearContents preserves the existing cell style (s="1") but drops the
# value/type, matching Excel's Delete-key behaviour on a formatted cell.
$ws.Range("B2:B9").ClearContents()
$ws.Range("D2:D9").ClearContents()
$ws.Range("F2:F9").ClearContents()

# C2 ("기초창의공학설계") was carrying the plain default style; give it the
# same Korean-charset font formatting already used on A2 (and every other
# styled cell in the table) by copying formats only.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Selection moved from F9 to E7 by the time the file was saved again.
$ws.Range("E7").Select() | Out-Null
